$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Consolidate the "Avatar of Discord" card rows (A2:A8) into a single
# Python-tuple-like string in A2.
$avatar = '(''Avatar of Discord'', [''{B/R}{B/R}{B/R}'', ''Creature — Avatar'', ''({B/R} can be paid with either {B} or {R}.)'', ''Flying'', ''When Avatar of Discord enters the battlefield, sacrifice it unless you discard two cards.'', ''5/3''])'

# Consolidate the "Azorius Guildmage" card rows (A9:A14) into a single
# Python-tuple-like string in A3.
$guildmage = '(''Azorius Guildmage'', [''{W/U}{W/U}'', ''Creature — Vedalken Wizard'', ''{2}{W}: Tap target creature.'', ''{2}{U}: Counter target activated ability. (Mana abilities can’t be targeted.)'', ''2/2''])'

$ws.Range("A2").Value = $avatar
$ws.Range("A3").Value = $guildmage

# Remove the now-redundant rows that held the individual card fields
# (old rows 4 through 14), shifting everything up.
$ws.Rows("4:14").Delete()
